$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list values (prices and 1h volume %) per the scraped data refresh,
# including the ShibaInu/WrappedEther row swap (rows 16 and 17).
# Values are written with a leading apostrophe so Excel keeps numeric-looking
# strings (prices, percentages) as text, matching the inlineStr cells in the
# source sheet; the Style reset drops the resulting quote-prefix flag so the
# cell formatting/style index is left untouched.

$ws.Range("D2").Value = '''61.074.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.07%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.397.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -1.63%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.04%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''573.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.92%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''142.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.32%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''3.398.48'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -1.63%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +0.06%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.476'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.77%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''7.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -1.05%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -2.75%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.397'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +1.27%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''3.981.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -1.52%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = '''  +2.03%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''28.08'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -2.44%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = '''WrappedEther'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '''3.407.73'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -2.49%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = '''ShibaInu'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = '''https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '''0.0000171'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.90%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''61.099.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -1.11%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''6.13'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -4.16%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''13.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -3.67%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''8.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -5.26%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''383.60'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -5.58%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.558'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -1.75%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''74.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +0.11%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.06%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  -5.16%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''3.538.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = '''  -1.64%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -0.32%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''7.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.60%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''8.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -3.22%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -1.80%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''1.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -4.30%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D35").Value = '''23.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -2.30%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''7.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -1.16%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''167.89'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +0.57%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''3.429.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -1.40%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''5.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -3.26%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''1.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -5.90%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.0773'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -2.92%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''27.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = '''  -2.73%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  +0.08%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''4.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -2.30%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  -4.35%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''  -1.81%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''2.479.99'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -5.01%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -2.45%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''23.06'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -0.79%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  +0.80%  '
$ws.Range("E51").Style = "Normal"
